# Updates market-price derived columns (H:N) across the Tonberry Profits
# leve-profit sheets, as produced by the scheduled price-update runner.
# Each entry: SheetName, Row, then a hashtable of Column -> Value (or $null
# to clear/remove the cell entirely, matching rows that lost a column).

$wb = $excel.ActiveWorkbook

$edits = @(
    @{ Sheet = "ALC"; Row = 6;   Cols = @{ H = 1948.6875;  J = 2565.75;    L = 7697.25;        N = -7921.25 } },
    @{ Sheet = "ALC"; Row = 51;  Cols = @{ H = 6000;       I = 5000;       J = 7000;   K = 5000; L = 7000;   M = -4516;             N = -7968 } },
    @{ Sheet = "ALC"; Row = 64;  Cols = @{ H = 0;          I = 0;          J = 0;      K = 0;    L = 0;      M = $null;             N = $null } },
    @{ Sheet = "ALC"; Row = 67;  Cols = @{ H = 0;          I = 0;          J = 0;      K = 0;    L = 0;      M = $null;             N = $null } },
    @{ Sheet = "ALC"; Row = 70;  Cols = @{ H = 1320;       I = 1000;       J = 1480;   K = 3000; L = 4440;   M = -2730;             N = -4980 } },
    @{ Sheet = "ALC"; Row = 73;  Cols = @{ H = 1320;       I = 1000;       J = 1480;   K = 3000; L = 4440;   M = -2064;             N = -6312 } },
    @{ Sheet = "ALC"; Row = 129; Cols = @{ H = 900.55554;  J = 875.58826;  L = 2626.76478;     N = -12626.76478 } },
    @{ Sheet = "ALC"; Row = 131; Cols = @{ H = 3402.5833;  I = 927.8;      J = 5170.2856; K = 2783.4; L = 15510.8568; M = 2256.6;  N = -25590.8568 } },
    @{ Sheet = "ALC"; Row = 132; Cols = @{ H = 1170.4166;  I = 1170.4166;  K = 3511.2498; M = -981.2498000000001 } },
    @{ Sheet = "ALC"; Row = 137; Cols = @{ H = 1473.7084;  I = 1255.7142;  J = 2999.6667; K = 3767.1426; L = 8999.000100000001; M = -1217.1426; N = -14099.0001 } },

    @{ Sheet = "ARM"; Row = 4;   Cols = @{ H = 173.33333;  I = 173.33333;  K = 173.33333; M = -58.33332999999999 } },
    @{ Sheet = "ARM"; Row = 5;   Cols = @{ H = 173.33333;  I = 173.33333;  K = 173.33333; M = -61.33332999999999 } },
    @{ Sheet = "ARM"; Row = 45;  Cols = @{ H = 1561.5;     I = 1128.8572;  K = 1128.8572; M = -751.8571999999999 } },
    @{ Sheet = "ARM"; Row = 61;  Cols = @{ H = 3179.1155;  I = 2415.2917;  K = 2415.2917; M = -2203.2917 } },
    @{ Sheet = "ARM"; Row = 136; Cols = @{ H = 3179.1155;  I = 2415.2917;  K = 7245.875100000001; M = -4695.875100000001 } },

    @{ Sheet = "BSM"; Row = 22;  Cols = @{ H = 599.6;      I = 599.5;      K = 599.5;  M = -426.5 } },

    @{ Sheet = "CRP"; Row = 7;   Cols = @{ I = 119.42857;  J = 71.5;       K = 119.42857; L = 71.5; M = -6.428569999999993; N = -297.5 } },
    @{ Sheet = "CRP"; Row = 31;  Cols = @{ H = 2483.1843;  I = 1690.5555;  J = 4428.727; K = 1690.5555; L = 4428.727; M = -1395.5555; N = -5018.727 } },
    @{ Sheet = "CRP"; Row = 34;  Cols = @{ H = 2483.1843;  I = 1690.5555;  J = 4428.727; K = 1690.5555; L = 4428.727; M = -1488.5555; N = -4832.727 } },
    @{ Sheet = "CRP"; Row = 56;  Cols = @{ H = 12000;      I = 12000;      K = 12000;  M = -11155 } },
    @{ Sheet = "CRP"; Row = 59;  Cols = @{ H = 20000;      J = 20000;      L = 20000;  N = -22290 } },
    @{ Sheet = "CRP"; Row = 62;  Cols = @{ H = 2975;       I = 0;          K = 0;      M = $null } },
    @{ Sheet = "CRP"; Row = 65;  Cols = @{ H = 2975;       I = 0;          K = 0;      M = $null } },
    @{ Sheet = "CRP"; Row = 122; Cols = @{ H = 1686.5416;  I = 1774.6;     J = 1539.7778; K = 5323.799999999999; L = 4619.3334; M = -2873.799999999999; N = -9519.3334 } },
    @{ Sheet = "CRP"; Row = 132; Cols = @{ H = 1734.3024;  I = 960.6;      J = 3519.7693; K = 2881.8; L = 10559.3079; M = -351.8000000000002; N = -15619.3079 } },

    @{ Sheet = "CUL"; Row = 131; Cols = @{ H = 9668.989;   J = 10432.397;  L = 31297.191; N = -41377.19100000001 } },

    @{ Sheet = "GSM"; Row = 70;  Cols = @{ H = 5133.3335;  I = 7000;       K = 7000;   M = -6730 } },
    @{ Sheet = "GSM"; Row = 73;  Cols = @{ H = 5133.3335;  I = 7000;       K = 7000;   M = -6064 } },

    @{ Sheet = "LTW"; Row = 40;  Cols = @{ H = 4982.6665;  I = 1598.1428;  K = 1598.1428; M = -1462.1428 } },
    @{ Sheet = "LTW"; Row = 61;  Cols = @{ H = 4166.6665;  I = 2500;       J = 5000;   K = 2500; L = 5000; M = -2298; N = -5404 } },
    @{ Sheet = "LTW"; Row = 74;  Cols = @{ H = 0;          J = 0;          L = 0;      N = $null } },
    @{ Sheet = "LTW"; Row = 75;  Cols = @{ H = 0;          J = 0;          L = 0;      N = $null } },
    @{ Sheet = "LTW"; Row = 77;  Cols = @{ H = 0;          J = 0;          L = 0;      N = $null } },
    @{ Sheet = "LTW"; Row = 78;  Cols = @{ H = 0;          J = 0;          L = 0;      N = $null } },
    @{ Sheet = "LTW"; Row = 113; Cols = @{ H = 4166.6665;  I = 2500;       J = 5000;   K = 2500; L = 5000; M = -330; N = -9340 } },
    @{ Sheet = "LTW"; Row = 122; Cols = @{ H = 5816.3887;  I = 4641.25;    K = 13923.75; M = -11473.75 } },
    @{ Sheet = "LTW"; Row = 132; Cols = @{ H = 1949.2222;  I = 1768.5385;  J = 2117;   K = 5305.6155; L = 6351; M = -2775.6155; N = -11411 } },
    @{ Sheet = "LTW"; Row = 133; Cols = @{ H = 0;          J = 0;          L = 0;      N = $null } },
    @{ Sheet = "LTW"; Row = 134; Cols = @{ H = 0;          J = 0;          L = 0;      N = $null } },
    @{ Sheet = "LTW"; Row = 135; Cols = @{ H = 0;          J = 0;          L = 0;      N = $null } },
    @{ Sheet = "LTW"; Row = 136; Cols = @{ H = 4121.6665;  I = 2542;       J = 5250;   K = 7626; L = 15750; M = -5076; N = -20850 } },
    @{ Sheet = "LTW"; Row = 137; Cols = @{ H = 0;          J = 0;          L = 0;      N = $null } },

    @{ Sheet = "WVR"; Row = 113; Cols = @{ H = 1435.8572;  I = 1025.25;    K = 3075.75; M = -905.75 } },
    @{ Sheet = "WVR"; Row = 132; Cols = @{ H = 1817.7142;  I = 1482.8422;  J = 4999;   K = 4448.5266; L = 14997; M = -1918.5266; N = -20057 } },
    @{ Sheet = "WVR"; Row = 136; Cols = @{ H = 16341317;   J = 2783;       L = 8349;   N = -13449 } },
    @{ Sheet = "WVR"; Row = 137; Cols = @{ H = 0;          J = 0;          L = 0;      N = $null } },
    @{ Sheet = "WVR"; Row = 139; Cols = @{ H = 59943.08;   J = 59943.08;   L = 59943.08; N = -70223.08 } },
    @{ Sheet = "WVR"; Row = 141; Cols = @{ H = 74539.89999999999; J = 74539.89999999999; L = 74539.89999999999; N = -84899.89999999999 } }
)

foreach ($edit in $edits) {
    $ws = $wb.Worksheets.Item($edit.Sheet)
    foreach ($col in $edit.Cols.Keys) {
        $addr = "$col$($edit.Row)"
        $val = $edit.Cols[$col]
        if ($null -eq $val) {
            $ws.Range($addr).ClearContents()
        } else {
            $ws.Range($addr).Value = $val
        }
    }
}
